$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ingredient quantities to standard amounts for 4 servings (4x the original values)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0.5
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 0.5
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 400
$ws.Range("C12").Value = 400
$ws.Range("C13").Value = 4
$ws.Range("C14").Value = 400

# Update the active selection
$ws.Range("E10").Select()
